# Add team record (Wins/Losses/Ties) columns to the player data sheet.
# New columns: AD = Wins, AE = Losses, AF = Ties.
# Header row (row 1) gets the same style as the existing header cells.
# Every data row (2-51) gets the team's W/L/T record: 74 / 88 / 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 51

# --- Header row: copy formatting from an existing header cell (A1) so the
# new headers pick up the bold font / border / centered style (s="1"),
# then overwrite with the real header text.
$ws.Range("A1").Copy($ws.Range("AD1")) | Out-Null
$ws.Range("A1").Copy($ws.Range("AE1")) | Out-Null
$ws.Range("A1").Copy($ws.Range("AF1")) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows: write the constant team record to every player row.
$rowCount = $lastRow - 1

$wins = New-Object 'object[,]' $rowCount,1
$losses = New-Object 'object[,]' $rowCount,1
$ties = New-Object 'object[,]' $rowCount,1

for ($i = 0; $i -lt $rowCount; $i++) {
    $wins[$i,0] = 74
    $losses[$i,0] = 88
    $ties[$i,0] = 0
}

$ws.Range("AD2:AD$lastRow").Value = $wins
$ws.Range("AE2:AE$lastRow").Value = $losses
$ws.Range("AF2:AF$lastRow").Value = $ties
